# Add the "bat crowd" scene-quest event as a new row 40 in Sheet1, shifting
# the existing rows 40-47 down to 41-48, expanding the table accordingly,
# and bumping the level of the last event (Id 42040008) from 1 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a fresh row above current row 40 (keeps formatting of the row
#    that is being pushed down, which is what Excel normally does).
$ws.Rows.Item(40).Insert()

# 2. Fill in the data for the new "batcrowd" event on row 40.
$ws.Cells.Item(40, 1).Value = 42010030          # Id
$ws.Cells.Item(40, 2).Value = "蝙蝠群"           # Name
$ws.Cells.Item(40, 3).Value = 1                  # Type
$ws.Cells.Item(40, 4).Value = 0                  # Level
$ws.Cells.Item(40, 5).Value = 2                  # Danger
$ws.Cells.Item(40, 6).Value = "batcrowd"         # Ename
$ws.Cells.Item(40, 7).Value = "batcrowd"         # Figue
$ws.Cells.Item(40, 8).Value = "batcrowd"         # Script
$ws.Cells.Item(40, 17).Value = 43000035          # Q EnemyId
$ws.Cells.Item(40, 18).Value = "mini"            # R BattleMap
$ws.Cells.Item(40, 21).Value = 100               # U RewardGold
$ws.Cells.Item(40, 22).Value = 150               # V RewardFood
$ws.Cells.Item(40, 25).Value = 100               # Y RewardExp
$ws.Cells.Item(40, 34).Value = 100               # AH PunishHealth
$ws.Cells.Item(40, 35).Value = 150               # AI PunishMental

# Catalog (J) needs to be "战斗" like the row that used to be row 40 (now
# row 41); copy its formatting (fill color) over, then set the text.
$ws.Range("J41").Copy()
$ws.Range("J40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(40, 10).Value = "战斗"

# 3. The last row (old row 47, "Id" 42040008, now row 48) had its Level
#    bumped from 1 to 3 as part of this change.
$ws.Cells.Item(48, 3).Value = 3

# 4. Expand table "表3" to cover the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:AN48"))

# 5. Restore the view: scroll so row 16 is the first visible row under the
#    frozen header, and select E38 like in the edited workbook.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E38").Select() | Out-Null
